# Auto-generated: apply scheduled market-price/profit updates across all sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4508.9473
$ws.Range("I86").Value = 5002.2856
$ws.Range("J86").Value = 4221.1665
$ws.Range("K86").Value = 5002.2856
$ws.Range("L86").Value = 4221.1665
$ws.Range("M86").Value = -3879.2856
$ws.Range("N86").Value = -6467.1665
$ws.Range("H89").Value = 4508.9473
$ws.Range("I89").Value = 5002.2856
$ws.Range("J89").Value = 4221.1665
$ws.Range("K89").Value = 25011.428
$ws.Range("L89").Value = 21105.8325
$ws.Range("M89").Value = -19395.428
$ws.Range("N89").Value = -32337.8325
$ws.Range("H111").Value = 3809
$ws.Range("I111").Value = 2009.8
$ws.Range("J111").Value = 6058
$ws.Range("K111").Value = 6029.4
$ws.Range("L111").Value = 18174
$ws.Range("M111").Value = -2962.4
$ws.Range("N111").Value = -24308
$ws.Range("H132").Value = 5749.6665
$ws.Range("I132").Value = 4999.636
$ws.Range("K132").Value = 14998.908
$ws.Range("M132").Value = -12468.908
$ws.Range("H137").Value = 3501.9697
$ws.Range("I137").Value = 917.4375
$ws.Range("K137").Value = 2752.3125
$ws.Range("M137").Value = -202.3125
$ws.Range("H138").Value = 3265.652
$ws.Range("J138").Value = 3261.3157
$ws.Range("L138").Value = 9783.947100000001
$ws.Range("N138").Value = -20063.9471

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1860.0358
$ws.Range("I45").Value = 1036
$ws.Range("K45").Value = 1036
$ws.Range("M45").Value = -659
$ws.Range("H74").Value = 7582166
$ws.Range("I74").Value = 9616589
$ws.Range("J74").Value = 25740
$ws.Range("K74").Value = 9616589
$ws.Range("L74").Value = 25740
$ws.Range("M74").Value = -9615715
$ws.Range("N74").Value = -27488
$ws.Range("H77").Value = 7582166
$ws.Range("I77").Value = 9616589
$ws.Range("J77").Value = 25740
$ws.Range("K77").Value = 48082945
$ws.Range("L77").Value = 128700
$ws.Range("M77").Value = -48078577
$ws.Range("N77").Value = -137436
$ws.Range("H135").Value = 111993
$ws.Range("J135").Value = 111993
$ws.Range("L135").Value = 111993
$ws.Range("N135").Value = -122133

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 10065
$ws.Range("J49").Value = 10065
$ws.Range("L49").Value = 10065
$ws.Range("N49").Value = -10543
$ws.Range("H87").Value = 131000
$ws.Range("J87").Value = 131000
$ws.Range("L87").Value = 131000
$ws.Range("N87").Value = -133496
$ws.Range("H90").Value = 131000
$ws.Range("J90").Value = 131000
$ws.Range("L90").Value = 393000
$ws.Range("N90").Value = -405480
$ws.Range("H105").Value = 1762.5454
$ws.Range("I105").Value = 1266.4445
$ws.Range("K105").Value = 1266.4445
$ws.Range("M105").Value = 480.5554999999999
$ws.Range("H138").Value = 274387.5
$ws.Range("J138").Value = 274387.5
$ws.Range("L138").Value = 274387.5
$ws.Range("N138").Value = -284667.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H135").Value = 85000
$ws.Range("J135").Value = 85000
$ws.Range("L135").Value = 85000
$ws.Range("N135").Value = -95140
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H137").Value = 48166.332
$ws.Range("I137").Value = 12500
$ws.Range("K137").Value = 12500
$ws.Range("M137").Value = -7400
$ws.Range("H138").Value = 86000
$ws.Range("J138").Value = 86000
$ws.Range("L138").Value = 86000
$ws.Range("N138").Value = -96280

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 140.83333
$ws.Range("I10").Value = 148.6
$ws.Range("K10").Value = 445.8
$ws.Range("M10").Value = -306.8
$ws.Range("H55").Value = 5766.909
$ws.Range("I55").Value = 5990.6665
$ws.Range("J55").Value = 5612
$ws.Range("K55").Value = 17971.9995
$ws.Range("L55").Value = 16836
$ws.Range("M55").Value = -17794.9995
$ws.Range("N55").Value = -17190
$ws.Range("H75").Value = 523299.66
$ws.Range("J75").Value = 523299.66
$ws.Range("L75").Value = 1569898.98
$ws.Range("N75").Value = -1571894.98
$ws.Range("H78").Value = 523299.66
$ws.Range("J78").Value = 523299.66
$ws.Range("L78").Value = 4709696.939999999
$ws.Range("N78").Value = -4719680.939999999
$ws.Range("H113").Value = 1758.1666
$ws.Range("I113").Value = 1850
$ws.Range("J113").Value = 1746.6875
$ws.Range("K113").Value = 5550
$ws.Range("L113").Value = 5240.0625
$ws.Range("M113").Value = -3380
$ws.Range("N113").Value = -9580.0625
$ws.Range("H132").Value = 2120.7646
$ws.Range("I132").Value = 2110.818
$ws.Range("K132").Value = 18997.362
$ws.Range("M132").Value = -16467.362
$ws.Range("H141").Value = 112582.46
$ws.Range("I141").Value = 128846.336
$ws.Range("K141").Value = 386539.008
$ws.Range("M141").Value = -381359.008

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1276.4231
$ws.Range("I97").Value = 1410.0526
$ws.Range("K97").Value = 1410.0526
$ws.Range("M97").Value = -914.0526
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4060.9443
$ws.Range("I40").Value = 3121.6667
$ws.Range("K40").Value = 3121.6667
$ws.Range("M40").Value = -2985.6667
$ws.Range("H74").Value = 81249.375
$ws.Range("J74").Value = 105800
$ws.Range("L74").Value = 105800
$ws.Range("N74").Value = -107796
$ws.Range("H77").Value = 81249.375
$ws.Range("J77").Value = 105800
$ws.Range("L77").Value = 317400
$ws.Range("N77").Value = -327384
$ws.Range("H82").Value = 2279.2222
$ws.Range("I82").Value = 1403.5
$ws.Range("J82").Value = 2979.8
$ws.Range("K82").Value = 1403.5
$ws.Range("L82").Value = 2979.8
$ws.Range("M82").Value = -1042.5
$ws.Range("N82").Value = -3701.8
$ws.Range("H85").Value = 2279.2222
$ws.Range("I85").Value = 1403.5
$ws.Range("J85").Value = 2979.8
$ws.Range("K85").Value = 1403.5
$ws.Range("L85").Value = 2979.8
$ws.Range("M85").Value = -155.5
$ws.Range("N85").Value = -5475.8
$ws.Range("H93").Value = 66667604
$ws.Range("I93").Value = 71429360
$ws.Range("K93").Value = 71429360
$ws.Range("M93").Value = -71428112

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 28578114
$ws.Range("I62").Value = 6800
$ws.Range("J62").Value = 33340000
$ws.Range("K62").Value = 6800
$ws.Range("L62").Value = 33340000
$ws.Range("M62").Value = -6176
$ws.Range("N62").Value = -33341248
$ws.Range("H65").Value = 28578114
$ws.Range("I65").Value = 6800
$ws.Range("J65").Value = 33340000
$ws.Range("K65").Value = 34000
$ws.Range("L65").Value = 166700000
$ws.Range("M65").Value = -30880
$ws.Range("N65").Value = -166706240
$ws.Range("H93").Value = 100000
$ws.Range("J93").Value = 100000
$ws.Range("L93").Value = 100000
$ws.Range("N93").Value = -104992
$ws.Range("H126").Value = 4617.4287
$ws.Range("I126").Value = 1998.0588
$ws.Range("K126").Value = 5994.1764
$ws.Range("M126").Value = -3524.1764
